$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 6652
$ws.Range("E2").Value = 286
$ws.Range("F2").Value = 286
$ws.Range("G2").Value = 490
$ws.Range("H2").Value = 389
$ws.Range("I2").Value = 385
$ws.Range("J2").Value = 4
$ws.Range("K2").Value = 8908
$ws.Range("L2").Value = 2139
$ws.Range("M2").Value = 6768
$ws.Range("N2").Value = 6720
$ws.Range("O2").Value = 48
$ws.Range("P2").Value = 115
$ws.Range("Q2").Value = 171
$ws.Range("R2").Value = 214
$ws.Range("S2").Value = -439
$ws.Range("T2").Value = 187
$ws.Range("U2").Value = -16
$ws.Range("V2").Value = 406
$ws.Range("W2").Value = 4.3
$ws.Range("X2").Value = 5.85
$ws.Range("Y2").Value = 5.86
$ws.Range("Z2").Value = 4.35
$ws.Range("AA2").Value = 31.61
$ws.Range("AB2").Value = 5315.09
$ws.Range("AC2").Value = 1673
$ws.Range("AD2").Value = 11.71
$ws.Range("AE2").Value = 29220
$ws.Range("AF2").Value = 0.67
$ws.Range("AG2").Value = 150
$ws.Range("AH2").Value = 0.77
$ws.Range("AI2").Value = 9.04
$ws.Range("AJ2").Value = 16800000

$ws.Range("D3").Value = 6659
$ws.Range("E3").Value = 259
$ws.Range("F3").Value = 259
$ws.Range("G3").Value = 548
$ws.Range("H3").Value = 434
$ws.Range("I3").Value = 430
$ws.Range("J3").Value = 4
$ws.Range("K3").Value = 9184
$ws.Range("L3").Value = 2010
$ws.Range("M3").Value = 7174
$ws.Range("N3").Value = 7122
$ws.Range("O3").Value = 52
$ws.Range("P3").Value = 115
$ws.Range("Q3").Value = 342
$ws.Range("R3").Value = 8
$ws.Range("S3").Value = -60
$ws.Range("T3").Value = 41
$ws.Range("U3").Value = 301
$ws.Range("V3").Value = 385
$ws.Range("W3").Value = 3.89
$ws.Range("X3").Value = 6.52
$ws.Range("Y3").Value = 6.21
$ws.Range("Z3").Value = 4.8
$ws.Range("AA3").Value = 28.03
$ws.Range("AB3").Value = 5690.66
$ws.Range("AC3").Value = 1868
$ws.Range("AD3").Value = 9.05
$ws.Range("AE3").Value = 30968
$ws.Range("AF3").Value = 0.55
$ws.Range("AG3").Value = 175
$ws.Range("AH3").Value = 1.04
$ws.Range("AI3").Value = 9.41
$ws.Range("AJ3").Value = 19308690

$ws.Range("D4").Value = 6490
$ws.Range("E4").Value = 197
$ws.Range("F4").Value = 197
$ws.Range("G4").Value = 521
$ws.Range("H4").Value = 423
$ws.Range("I4").Value = 420
$ws.Range("J4").Value = 3
$ws.Range("K4").Value = 9398
$ws.Range("L4").Value = 1893
$ws.Range("M4").Value = 7505
$ws.Range("N4").Value = 7485
$ws.Range("O4").Value = 19
$ws.Range("P4").Value = 115
$ws.Range("Q4").Value = 250
$ws.Range("R4").Value = -144
$ws.Range("S4").Value = -222
$ws.Range("T4").Value = 233
$ws.Range("U4").Value = 17
$ws.Range("V4").Value = 272
$ws.Range("W4").Value = 3.03
$ws.Range("X4").Value = 6.52
$ws.Range("Y4").Value = 5.75
$ws.Range("Z4").Value = 4.56
$ws.Range("AA4").Value = 25.22
$ws.Range("AB4").Value = 6168.23
$ws.Range("AC4").Value = 1828
$ws.Range("AD4").Value = 8.29
$ws.Range("AE4").Value = 32548
$ws.Range("AF4").Value = 0.47
$ws.Range("AG4").Value = 175
$ws.Range("AH4").Value = 1.16
$ws.Range("AI4").Value = 9.62
$ws.Range("AJ4").Value = 19308690

$ws.Range("D5").Value = 6661
$ws.Range("E5").Value = 115
$ws.Range("F5").Value = 115
$ws.Range("G5").Value = 374
$ws.Range("H5").Value = 314
$ws.Range("I5").Value = 313
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 9538
$ws.Range("L5").Value = 1772
$ws.Range("M5").Value = 7766
$ws.Range("N5").Value = 7746
$ws.Range("O5").Value = 20
$ws.Range("P5").Value = 115
$ws.Range("Q5").Value = 293
$ws.Range("R5").Value = -71
$ws.Range("S5").Value = -95
$ws.Range("T5").Value = 154
$ws.Range("U5").Value = 139
$ws.Range("V5").Value = 265
$ws.Range("W5").Value = 1.72
$ws.Range("X5").Value = 4.72
$ws.Range("Y5").Value = 4.11
$ws.Range("Z5").Value = 3.32
$ws.Range("AA5").Value = 22.82
$ws.Range("AB5").Value = 6423.49
$ws.Range("AC5").Value = 1362
$ws.Range("AD5").Value = 9.25
$ws.Range("AE5").Value = 33963
$ws.Range("AF5").Value = 0.37
$ws.Range("AG5").Value = 175
$ws.Range("AH5").Value = 1.39
$ws.Range("AI5").Value = 12.8
$ws.Range("AJ5").Value = 19308690

$ws.Range("D6").Value = 6517
$ws.Range("E6").Value = 114
$ws.Range("F6").Value = 114
$ws.Range("G6").Value = 524
$ws.Range("H6").Value = 435
$ws.Range("I6").Value = 434
$ws.Range("K6").Value = 9889
$ws.Range("L6").Value = 1731
$ws.Range("M6").Value = 8158
$ws.Range("N6").Value = 8136
$ws.Range("P6").Value = 115
$ws.Range("Q6").Value = 251
$ws.Range("R6").Value = -264
$ws.Range("S6").Value = -59
$ws.Range("T6").Value = 335
$ws.Range("U6").Value = -84
$ws.Range("V6").Value = 246
$ws.Range("W6").Value = 1.75
$ws.Range("X6").Value = 6.68
$ws.Range("Y6").Value = 5.47
$ws.Range("Z6").Value = 4.48
$ws.Range("AA6").Value = 21.22
$ws.Range("AB6").Value = 6779.35
$ws.Range("AC6").Value = 1888
$ws.Range("AD6").Value = 6.36
$ws.Range("AE6").Value = 35675
$ws.Range("AF6").Value = 0.34
$ws.Range("AG6").Value = 175
$ws.Range("AH6").Value = 1.46
$ws.Range("AI6").Value = 9.24
$ws.Range("AJ6").Value = 19308690

$ws.Range("D7:AJ7").ClearContents()
$ws.Range("D8:AJ8").ClearContents()
$ws.Range("D9:AJ9").ClearContents()
